$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Direct writes for cells whose text will not be misinterpreted as numbers/dates
$ws.Range('D2').Value = '29.517.74'
$ws.Range('E2').Value = '  -1.35%  '
$ws.Range('D3').Value = '1.853.54'
$ws.Range('E3').Value = '  -0.49%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('E5').Value = '  -0.92%  '
$ws.Range('E6').Value = '  -1.97%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -0.36%  '
$ws.Range('E9').Value = '  -0.51%  '
$ws.Range('E10').Value = '  -1.39%  '
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').Value = '1.882.80'
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('E13').Value = '  -0.49%  '
$ws.Range('E14').Value = '  -1.28%  '
$ws.Range('E15').Value = '  -0.40%  '
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('E17').Value = '  +2.15%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '29.546.77'
$ws.Range('E18').Value = '  -1.28%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('E19').Value = '  -1.08%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('E20').Value = '  -1.28%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('E22').Value = '  -2.60%  '
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E24').Value = '  -2.11%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('E25').Value = '  -2.21%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('E26').Value = '  -1.51%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('E27').Value = '  -1.41%  '
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('E28').Value = '  -1.32%  '
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('E29').Value = '  -3.68%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('E30').Value = '  -3.06%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('E31').Value = '  -0.97%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('E32').Value = '  -1.39%  '
$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('E34').Value = '  -0.63%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('E35').Value = '  -0.84%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('E36').Value = '  -0.99%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('E37').Value = '  -0.56%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '1.240.58'
$ws.Range('E38').Value = '  +1.97%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('E40').Value = '  -1.45%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E41').Value = '  -4.34%  '
$ws.Range('B42').Value = 'RocketPoolETH'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D42').Value = '2.042.77'
$ws.Range('E42').Value = '  +0.27%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('E44').Value = '  +0.24%  '
$ws.Range('E45').Value = '  -0.72%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('E46').Value = '  -4.54%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('E47').Value = '  +0.92%  '
$ws.Range('B48').Value = 'TheSandbox'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('E48').Value = '  -1.00%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E49').Value = '  -1.01%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('E50').Value = '  +0.88%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('E51').Value = '  -0.64%  '

# Numeric-looking price strings must be forced to stay as literal text
# (matching the source data, which stores prices like "0.6330" or "1.000"
# as text, not numbers). Round-trip them through a scratch cell that has
# been explicitly formatted as Text, then paste-special the VALUE only back
# into the target cell so the target cell keeps its original (default) style.
$scratch = $ws.Range('Z1')
$scratch.NumberFormat = "@"
$scratch.Value = '0.9981'
$scratch.Copy()
$ws.Range('D4').PasteSpecial(-4163)
$scratch.Value = '242.12'
$scratch.Copy()
$ws.Range('D5').PasteSpecial(-4163)
$scratch.Value = '0.6330'
$scratch.Copy()
$ws.Range('D6').PasteSpecial(-4163)
$scratch.Value = '0.07539'
$scratch.Copy()
$ws.Range('D8').PasteSpecial(-4163)
$scratch.Value = '0.2978'
$scratch.Copy()
$ws.Range('D9').PasteSpecial(-4163)
$scratch.Value = '24.53'
$scratch.Copy()
$ws.Range('D10').PasteSpecial(-4163)
$scratch.Value = '0.07730'
$scratch.Copy()
$ws.Range('D11').PasteSpecial(-4163)
$scratch.Value = '0.6926'
$scratch.Copy()
$ws.Range('D13').PasteSpecial(-4163)
$scratch.Value = '5.018'
$scratch.Copy()
$ws.Range('D14').PasteSpecial(-4163)
$scratch.Value = '83.55'
$scratch.Copy()
$ws.Range('D15').PasteSpecial(-4163)
$scratch.Value = '0.000009853'
$scratch.Copy()
$ws.Range('D16').PasteSpecial(-4163)
$scratch.Value = '6.289'
$scratch.Copy()
$ws.Range('D17').PasteSpecial(-4163)
$scratch.Value = '233.89'
$scratch.Copy()
$ws.Range('D19').PasteSpecial(-4163)
$scratch.Value = '12.54'
$scratch.Copy()
$ws.Range('D20').PasteSpecial(-4163)
$scratch.Value = '0.9990'
$scratch.Copy()
$ws.Range('D21').PasteSpecial(-4163)
$scratch.Value = '7.672'
$scratch.Copy()
$ws.Range('D22').PasteSpecial(-4163)
$scratch.Value = '1.001'
$scratch.Copy()
$ws.Range('D23').PasteSpecial(-4163)
$scratch.Value = '155.48'
$scratch.Copy()
$ws.Range('D24').PasteSpecial(-4163)
$scratch.Value = '0.1396'
$scratch.Copy()
$ws.Range('D25').PasteSpecial(-4163)
$scratch.Value = '8.470'
$scratch.Copy()
$ws.Range('D26').PasteSpecial(-4163)
$scratch.Value = '17.74'
$scratch.Copy()
$ws.Range('D27').PasteSpecial(-4163)
$scratch.Value = '1.476'
$scratch.Copy()
$ws.Range('D28').PasteSpecial(-4163)
$scratch.Value = '0.05931'
$scratch.Copy()
$ws.Range('D29').PasteSpecial(-4163)
$scratch.Value = '1.253'
$scratch.Copy()
$ws.Range('D30').PasteSpecial(-4163)
$scratch.Value = '4.129'
$scratch.Copy()
$ws.Range('D31').PasteSpecial(-4163)
$scratch.Value = '4.051'
$scratch.Copy()
$ws.Range('D32').PasteSpecial(-4163)
$scratch.Value = '1.893'
$scratch.Copy()
$ws.Range('D33').PasteSpecial(-4163)
$scratch.Value = '1.168'
$scratch.Copy()
$ws.Range('D34').PasteSpecial(-4163)
$scratch.Value = '0.7230'
$scratch.Copy()
$ws.Range('D35').PasteSpecial(-4163)
$scratch.Value = '2.585'
$scratch.Copy()
$ws.Range('D36').PasteSpecial(-4163)
$scratch.Value = '2.802'
$scratch.Copy()
$ws.Range('D37').PasteSpecial(-4163)
$scratch.Value = '0.01795'
$scratch.Copy()
$ws.Range('D39').PasteSpecial(-4163)
$scratch.Value = '0.9068'
$scratch.Copy()
$ws.Range('D40').PasteSpecial(-4163)
$scratch.Value = '6.105'
$scratch.Copy()
$ws.Range('D41').PasteSpecial(-4163)
$scratch.Value = '0.9994'
$scratch.Copy()
$ws.Range('D43').PasteSpecial(-4163)
$scratch.Value = '67.58'
$scratch.Copy()
$ws.Range('D44').PasteSpecial(-4163)
$scratch.Value = '101.56'
$scratch.Copy()
$ws.Range('D45').PasteSpecial(-4163)
$scratch.Value = '7.415'
$scratch.Copy()
$ws.Range('D46').PasteSpecial(-4163)
$scratch.Value = '0.00000000120'
$scratch.Copy()
$ws.Range('D47').PasteSpecial(-4163)
$scratch.Value = '0.4047'
$scratch.Copy()
$ws.Range('D48').PasteSpecial(-4163)
$scratch.Value = '9.130'
$scratch.Copy()
$ws.Range('D49').PasteSpecial(-4163)
$scratch.Value = '1.709'
$scratch.Copy()
$ws.Range('D50').PasteSpecial(-4163)
$scratch.Value = '0.05757'
$scratch.Copy()
$ws.Range('D51').PasteSpecial(-4163)
$scratch.Clear()

